$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 72; this shifts existing rows 72-180 down to 73-181,
# preserving their data (matches the bulk of the diff which is just a row-index shift).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data point.
$ws.Range("A72").Value = 8
$ws.Range("B72").Value = "Terminal La Palmera de La Serena"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44671
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E72").Value = 4
$ws.Range("F72").Value = 100112037
$ws.Range("G72").Value = "Cebollín"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 1000
$ws.Range("K72").Value = 1100
$ws.Range("L72").Value = 1200
$ws.Range("M72").Value = 1150
$ws.Range("N72").Value = "$/paquete 6 unidades"
$ws.Range("O72").Value = "Provincia del Elquí"
$ws.Range("P72").Value = 192
$ws.Range("Q72").Value = 6
$ws.Range("R72").Value = "Hortaliza"
